$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 16: Rotary joint mount now has a supplier (Farnell) and a plain
#     (non right-aligned) numeric part number ---
$ws.Range("D16").Value = "Farnell"
$ws.Range("E16").Value = 1924856
$ws.Range("E16").Style = "Normal"

# --- Row 26: Lens swapped from Computar 8mm CS mount lens to the
#     2.8-12mm CS mount lens ---
$ws.Range("B26").Value = "2.8-12mm CS mount Lens"
$ws.Range("E26").Value = "T4Z2813CS-IR"

# --- Row 27: Ball head mount swapped from the SIOTI (Amazon link) mount
#     to the Ulanzi mount, supplier UK Digital, part UL-U30. The old
#     hyperlink cell no longer needs the Hyperlink style. ---
$ws.Range("B27").Value = "Ulanzi Mini Ball head mount"
$ws.Range("D27").Value = "UK Digital"
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").WrapText = $true
$ws.Range("E27").Value = "UL-U30"

# The Hyperlink built-in cell style is no longer used anywhere in the
# workbook now that the Amazon link text has been replaced, so remove it.
$wb.Styles.Item("Hyperlink").Delete()

# --- Row 28: Cap head bolt now has a supplier (ACCU group) and part
#     number (SSC-1/4-20-1-A2) instead of placeholder dashes ---
$ws.Range("D28").Value = "ACCU group"
$ws.Range("E28").Value = "SSC-1/4-20-1-A2"

# --- Restore the sheet view selection state ---
$ws.Range("B28").Select()
